# Updated cryptos list on Thu Aug  1 07:53:12 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with the latest
# figures pulled from coinranking.com, and fixes up a few rows whose
# rank changed (InjectiveProtocol/Bittensor and dogwifhat/FirstDigitalUSD
# swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price (D) and volume (E) values for rows with changed figures.
# Some "Price" values are plain decimal numbers (e.g. 569.85, 0.0620);
# force those cells to Text format first so Excel keeps the original
# textual representation (trailing zeros, no scientific notation) instead
# of silently converting them to floating point numbers.
$ws.Range("D2").Value = "64.386.86"
$ws.Range("E2").Value = "  -2.95%  "
$ws.Range("D3").Value = "3.172.41"
$ws.Range("E3").Value = "  -4.59%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.85"
$ws.Range("E5").Value = "  -2.83%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.42"
$ws.Range("E6").Value = "  -8.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  -6.09%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "3.172.76"
$ws.Range("E9").Value = "  -4.48%  "
$ws.Range("E10").Value = "  -4.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.77"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.386"
$ws.Range("E12").Value = "  -4.25%  "
$ws.Range("D13").Value = "3.724.44"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.129"
$ws.Range("E14").Value = "  -2.55%  "
$ws.Range("D15").Value = "64.429.98"
$ws.Range("E15").Value = "  -2.92%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.34"
$ws.Range("E16").Value = "  -4.03%  "
$ws.Range("E17").Value = "  -3.07%  "
$ws.Range("D18").Value = "3.173.19"
$ws.Range("E18").Value = "  -4.97%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "417.86"
$ws.Range("E19").Value = "  -2.88%  "
$ws.Range("E20").Value = "  -3.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.85"
$ws.Range("E21").Value = "  -3.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.06"
$ws.Range("E22").Value = "  -5.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.75"
$ws.Range("E24").Value = "  -3.26%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("E26").Value = "  -6.16%  "
$ws.Range("E27").Value = "  -7.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.87"
$ws.Range("E28").Value = "  -1.50%  "
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("E30").Value = "  -6.05%  "
$ws.Range("E31").Value = "  -3.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.03"
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.32"
$ws.Range("E34").Value = "  -4.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.14"
$ws.Range("E35").Value = "  -4.78%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "157.12"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("E37").Value = "  -6.66%  "
$ws.Range("D38").Value = "2.726.16"
$ws.Range("E38").Value = "  -6.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.70"
$ws.Range("E39").Value = "  -6.76%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "24.30"
$ws.Range("E40").Value = "  -9.03%  "
$ws.Range("E41").Value = "  -4.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.13"
$ws.Range("E43").Value = "  -7.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0620"
$ws.Range("E44").Value = "  -7.28%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.59"
$ws.Range("E45").Value = "  -6.93%  "
$ws.Range("E46").Value = "  -4.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0986"
$ws.Range("E51").Value = "  -6.09%  "

# Rows 47-50 changed order/content: some coins moved up in ranking,
# so refresh the Coin/Link/Price/Volume columns for those rows directly.
$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "293.55"
$ws.Range("E47").Value = "  -7.19%  "

$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.61"
$ws.Range("E48").Value = "  -7.73%  "

$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.01"
$ws.Range("E49").Value = "  -13.27%  "

$ws.Range("B50").Value = "FirstDigitalUSD"
$ws.Range("C50").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.01%  "
